$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "SA"

# Add new row 16 data, reusing the formatting from row 15 (A15 has a bold,
# centered, bordered style applied) for the new index cell A16.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A16").Value = 14

$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"

$ws.Range("C16").Value = 1.02888972242382
$ws.Range("D16").Value = 0.9467527183466591
$ws.Range("E16").Value = 1.006267620948166
$ws.Range("F16").Value = 0.9859838365847405
$ws.Range("G16").Value = 1.02888972242382
$ws.Range("H16").Value = 0.9467527183466591
$ws.Range("I16").Value = 1.010908137415083
$ws.Range("J16").Value = 0.9866228137006403
$ws.Range("K16").Value = 1.006267620948166
$ws.Range("L16").Value = 0.9641314439970358
$ws.Range("M16").Value = 1.02888972242382
$ws.Range("N16").Value = 0.9765101696474126
$ws.Range("O16").Value = 0.9919734745758463
$ws.Range("P16").Value = 0.9919779892955387
